$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in row 1 (D1/E1), splitting the old combined columns
$ws.Range("D1").Value = "No of classes alloted"
$ws.Range("E1").Value = "No of classes taken"

# Update the data value in E2 (was a percentage 0-100, now a plain count)
$ws.Range("E2").Value = 11

# Update the selected cell/range shown when the file is opened
$ws.Range("E2").Select()

# Replace the data validation on column E (was decimal between 0 and 100,
# now decimal greater than -1)
$eRange = $ws.Range("E2:E1048576")
$eRange.Validation.Delete()
$eRange.Validation.Add(2, 1, 5, "-1")
